# "updated legacy GSC export data"
#
# The "Chart" sheet holds a rolling 90-day window of GSC export data:
#   column A = date (stored as text, e.g. "2025-09-12")
#   column B = Non-HTTPS URLs (always 0 in this export)
#   column C = HTTPS URLs (count)
#
# The window rolled forward by one day: the oldest date (2025-09-12, row 2)
# drops off, every remaining row shifts up by one, and a new row for the
# newest date (2025-12-11) is appended at the bottom with its URL count.
#
# We replicate that by copying each row's A/C values from the row below it
# (walking top-down so we always read a row before it gets overwritten),
# then writing the new trailing date/value pair into the last row.

$wb = $excel.ActiveWorkbook
$chart = $wb.Worksheets.Item("Chart")

$firstDataRow = 2
$lastDataRow = 91

for ($r = $firstDataRow; $r -le ($lastDataRow - 1); $r++) {
    $srcRow = $r + 1

    # Read the date as displayed text (not .Value, which would hand back a
    # date serial once the destination cell is re-typed) and the URL count
    # as a plain number.
    $dateText = $chart.Cells.Item($srcRow, 1).Text
    $count = $chart.Cells.Item($srcRow, 3).Value2

    # Force the destination to stay text so the date string isn't
    # reinterpreted/converted into a date serial on write.
    $chart.Cells.Item($r, 1).NumberFormat = "@"
    $chart.Cells.Item($r, 1).Value = $dateText
    $chart.Cells.Item($r, 3).Value = $count
}

# Append the new newest day at the end of the window.
$chart.Cells.Item($lastDataRow, 1).NumberFormat = "@"
$chart.Cells.Item($lastDataRow, 1).Value = "2025-12-11"
$chart.Cells.Item($lastDataRow, 3).Value = 0
